# Generate Report for Handoff
#
# The f1934dc0-... row moves from "Handed back: in sync with en-US" to
# "Ready for handoff" in both locale sheets (and the Overview rollup),
# gets a refreshed "Latest Handoff Datetime" stamp, and picks up an
# Error Detail message explaining the handback file is stale. The
# Error Detail column is widened to fit the longer message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e601a323243a7e78fbb598716b28352803c4440e/e2e/f1934dc0-44fe-469d-81ce-6c5a14bf8cf2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db4eb3cd99593ee7e637336c414826b8cfa516b2/e2e/f1934dc0-44fe-469d-81ce-6c5a14bf8cf2.md."

# ---- Overview sheet: row 3 is the f1934dc0-... file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-26 16:48:58"

# ---- zh-cn sheet: row 3 is the f1934dc0-... file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("H3").Value = "2016-08-26 16:48:54"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.14

# ---- de-de sheet: row 3 is the f1934dc0-... file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("H3").Value = "2016-08-26 16:48:58"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.14
